# Update the Sources list on the "Sources (引用元一覧)" slide (slide 15).
# The list of source citations was re-ordered/rotated: each bullet's
# text is replaced with the text that should now appear in that position.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(15)
$shp = $s.Shapes.Item(4)
$tr = $shp.TextFrame.TextRange

$tr.Paragraphs(2).Runs(1).Text = "CNN Business: https://www.cnn.com/2026/02/05/tech/anthropic-opus-update-software-stocks"
$tr.Paragraphs(5).Runs(1).Text = "MIT Sloan Management Review: https://sloanreview.mit.edu/article/five-trends-in-ai-and-data-science-for-2026/"
$tr.Paragraphs(6).Runs(1).Text = "Tech Startups: https://techstartups.com/2026/02/05/top-startup-and-tech-funding-news-february-5-2025/"
$tr.Paragraphs(7).Runs(1).Text = "LLM Stats: https://llm-stats.com/llm-updates"
$tr.Paragraphs(8).Runs(1).Text = "National Law Review: https://natlawreview.com/article/2026-outlook-artificial-intelligence"
$tr.Paragraphs(9).Runs(1).Text = "MIT Technology Review: https://www.technologyreview.com/2026/01/05/1130662/whats-next-for-ai-in-2026/"
$tr.Paragraphs(10).Runs(1).Text = "TechCrunch: https://techcrunch.com/2026/02/05/openai-launches-new-agentic-coding-model-only-minutes-after-anthropic-drops-its-own/"
$tr.Paragraphs(11).Runs(1).Text = "Axios: https://www.axios.com/2026/02/06/amazon-microsoft-meta-ai-investment"
